{"js": "// The document contains a single table of 20 rows x 5 columns, each cell\n// holding a short arithmetic expression like \"58-29=\". The commit replaces\n// every one of the 100 expressions with a new one (same cell position,\n// same formatting) \u2014 this is the \"adc-sbb-within100\" (addition/subtraction\n// within 100) worksheet being regenerated with a fresh batch of problems.\n//\n// Old -> new text, in document (row-major) order: row 1 col 1..5, row 2\n// col 1..5, ... row 20 col 1..5. A couple of \"old\" strings repeat (e.g.\n// \"58+25=\" and \"38+7=\" each occur twice with different replacements), so\n// the substitution must be positional, not a global text search/replace.\nconst oldNewPairs = [\n  [\"58-29=\", \"37+29=\"], [\"42-37=\", \"51-7=\"], [\"91-49=\", \"29+68=\"], [\"25+7=\", \"70-55=\"], [\"68+19=\", \"61-25=\"],\n  [\"74+19=\", \"72-64=\"], [\"59+38=\", \"91-73=\"], [\"57+25=\", \"84+8=\"], [\"8+7=\", \"52-38=\"], [\"23+29=\", \"35+46=\"],\n  [\"58+25=\", \"17+36=\"], [\"12+29=\", \"35+47=\"], [\"13+78=\", \"3+18=\"], [\"84-67=\", \"93-28=\"], [\"38+29=\", \"57+28=\"],\n  [\"94-65=\", \"43-19=\"], [\"38+7=\", \"3+69=\"], [\"37+45=\", \"41-17=\"], [\"53-16=\", \"90-19=\"], [\"73-66=\", \"81-75=\"],\n  [\"49+24=\", \"90-65=\"], [\"53-17=\", \"17+36=\"], [\"13+18=\", \"48+16=\"], [\"83-27=\", \"76+9=\"], [\"12-8=\", \"42+9=\"],\n  [\"39+25=\", \"60-39=\"], [\"51-13=\", \"51-23=\"], [\"5+28=\", \"28+45=\"], [\"26+59=\", \"93-47=\"], [\"7+68=\", \"28+19=\"],\n  [\"9+24=\", \"9+85=\"], [\"70-44=\", \"30-17=\"], [\"40-31=\", \"60-37=\"], [\"55-26=\", \"85-29=\"], [\"49+17=\", \"82-24=\"],\n  [\"85+7=\", \"58+27=\"], [\"37+35=\", \"43+39=\"], [\"96-39=\", \"90-87=\"], [\"70-26=\", \"57+18=\"], [\"62-5=\", \"91-89=\"],\n  [\"90-41=\", \"37+19=\"], [\"58+15=\", \"68-29=\"], [\"65+28=\", \"28+46=\"], [\"33-24=\", \"71-22=\"], [\"3+29=\", \"94-66=\"],\n  [\"14+78=\", \"7+59=\"], [\"15+26=\", \"29+13=\"], [\"47+7=\", \"93-74=\"], [\"58+25=\", \"96-28=\"], [\"92-6=\", \"16+39=\"],\n  [\"54+17=\", \"37-8=\"], [\"47+9=\", \"92-45=\"], [\"97-39=\", \"81-52=\"], [\"95-16=\", \"50-49=\"], [\"67+26=\", \"29+18=\"],\n  [\"75+9=\", \"84-9=\"], [\"69+12=\", \"37+47=\"], [\"69+16=\", \"92-43=\"], [\"48+25=\", \"53-34=\"], [\"49+33=\", \"24+18=\"],\n  [\"67+4=\", \"43+28=\"], [\"68+13=\", \"40-11=\"], [\"59+19=\", \"66-59=\"], [\"17+45=\", \"71-6=\"], [\"74-38=\", \"48+45=\"],\n  [\"38+7=\", \"2+59=\"], [\"83-44=\", \"86+7=\"], [\"66-18=\", \"58+29=\"], [\"87+5=\", \"89+4=\"], [\"14+38=\", \"40-37=\"],\n  [\"12+49=\", \"76+8=\"], [\"34-9=\", \"56-29=\"], [\"8+84=\", \"15+17=\"], [\"92-83=\", \"31-6=\"], [\"39+49=\", \"36+17=\"],\n  [\"9+23=\", \"63+28=\"], [\"67+8=\", \"78+6=\"], [\"51-27=\", \"69+14=\"], [\"23+59=\", \"33-8=\"], [\"29+23=\", \"3+88=\"],\n  [\"19+68=\", \"88+8=\"], [\"95-39=\", \"19+56=\"], [\"18+79=\", \"33+58=\"], [\"2+9=\", \"37+54=\"], [\"49+6=\", \"22-14=\"],\n  [\"96-78=\", \"81-49=\"], [\"28+8=\", \"73-44=\"], [\"75-9=\", \"82-3=\"], [\"30-16=\", \"48+47=\"], [\"82-77=\", \"90-16=\"],\n  [\"33-14=\", \"65+9=\"], [\"46-17=\", \"6+19=\"], [\"91-86=\", \"3+49=\"], [\"57+5=\", \"88-69=\"], [\"15+68=\", \"76-39=\"],\n  [\"75-39=\", \"45+8=\"], [\"90-12=\", \"16+65=\"], [\"71-62=\", \"3+79=\"], [\"61-58=\", \"70-29=\"], [\"25-18=\", \"36+5=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst current = table.values;\n\nlet i = 0;\nconst updated = current.map((row) =>\n  row.map((cellText) => {\n    const [oldText, newText] = oldNewPairs[i];\n    i++;\n    // Sanity check the expected original text; if it doesn't line up\n    // (shouldn't happen for this document), keep the existing cell text\n    // untouched rather than corrupting unrelated content.\n    return cellText === oldText ? newText : cellText;\n  })\n);\n\ntable.values = updated;\nawait context.sync();\n", "ps1": "# The document body contains a single 20-row x 5-column table of short\n# arithmetic expressions (e.g. \"58-29=\"). The commit swaps every one of\n# the 100 expressions for a newly generated one, in place, keeping each\n# cell's run formatting (font/size) and paragraph formatting untouched.\n#\n# Mapping is positional (row-major: row 1 col 1..5, row 2 col 1..5, ...)\n# rather than a global text replace, because a few \"old\" strings repeat\n# in the grid with different replacements (e.g. \"58+25=\" occurs twice).\n$pairs = @(\n    @{ Old = \"58-29=\"; New = \"37+29=\" },\n    @{ Old = \"42-37=\"; New = \"51-7=\" },\n    @{ Old = \"91-49=\"; New = \"29+68=\" },\n    @{ Old = \"25+7=\"; New = \"70-55=\" },\n    @{ Old = \"68+19=\"; New = \"61-25=\" },\n    @{ Old = \"74+19=\"; New = \"72-64=\" },\n    @{ Old = \"59+38=\"; New = \"91-73=\" },\n    @{ Old = \"57+25=\"; New = \"84+8=\" },\n    @{ Old = \"8+7=\"; New = \"52-38=\" },\n    @{ Old = \"23+29=\"; New = \"35+46=\" },\n    @{ Old = \"58+25=\"; New = \"17+36=\" },\n    @{ Old = \"12+29=\"; New = \"35+47=\" },\n    @{ Old = \"13+78=\"; New = \"3+18=\" },\n    @{ Old = \"84-67=\"; New = \"93-28=\" },\n    @{ Old = \"38+29=\"; New = \"57+28=\" },\n    @{ Old = \"94-65=\"; New = \"43-19=\" },\n    @{ Old = \"38+7=\"; New = \"3+69=\" },\n    @{ Old = \"37+45=\"; New = \"41-17=\" },\n    @{ Old = \"53-16=\"; New = \"90-19=\" },\n    @{ Old = \"73-66=\"; New = \"81-75=\" },\n    @{ Old = \"49+24=\"; New = \"90-65=\" },\n    @{ Old = \"53-17=\"; New = \"17+36=\" },\n    @{ Old = \"13+18=\"; New = \"48+16=\" },\n    @{ Old = \"83-27=\"; New = \"76+9=\" },\n    @{ Old = \"12-8=\"; New = \"42+9=\" },\n    @{ Old = \"39+25=\"; New = \"60-39=\" },\n    @{ Old = \"51-13=\"; New = \"51-23=\" },\n    @{ Old = \"5+28=\"; New = \"28+45=\" },\n    @{ Old = \"26+59=\"; New = \"93-47=\" },\n    @{ Old = \"7+68=\"; New = \"28+19=\" },\n    @{ Old = \"9+24=\"; New = \"9+85=\" },\n    @{ Old = \"70-44=\"; New = \"30-17=\" },\n    @{ Old = \"40-31=\"; New = \"60-37=\" },\n    @{ Old = \"55-26=\"; New = \"85-29=\" },\n    @{ Old = \"49+17=\"; New = \"82-24=\" },\n    @{ Old = \"85+7=\"; New = \"58+27=\" },\n    @{ Old = \"37+35=\"; New = \"43+39=\" },\n    @{ Old = \"96-39=\"; New = \"90-87=\" },\n    @{ Old = \"70-26=\"; New = \"57+18=\" },\n    @{ Old = \"62-5=\"; New = \"91-89=\" },\n    @{ Old = \"90-41=\"; New = \"37+19=\" },\n    @{ Old = \"58+15=\"; New = \"68-29=\" },\n    @{ Old = \"65+28=\"; New = \"28+46=\" },\n    @{ Old = \"33-24=\"; New = \"71-22=\" },\n    @{ Old = \"3+29=\"; New = \"94-66=\" },\n    @{ Old = \"14+78=\"; New = \"7+59=\" },\n    @{ Old = \"15+26=\"; New = \"29+13=\" },\n    @{ Old = \"47+7=\"; New = \"93-74=\" },\n    @{ Old = \"58+25=\"; New = \"96-28=\" },\n    @{ Old = \"92-6=\"; New = \"16+39=\" },\n    @{ Old = \"54+17=\"; New = \"37-8=\" },\n    @{ Old = \"47+9=\"; New = \"92-45=\" },\n    @{ Old = \"97-39=\"; New = \"81-52=\" },\n    @{ Old = \"95-16=\"; New = \"50-49=\" },\n    @{ Old = \"67+26=\"; New = \"29+18=\" },\n    @{ Old = \"75+9=\"; New = \"84-9=\" },\n    @{ Old = \"69+12=\"; New = \"37+47=\" },\n    @{ Old = \"69+16=\"; New = \"92-43=\" },\n    @{ Old = \"48+25=\"; New = \"53-34=\" },\n    @{ Old = \"49+33=\"; New = \"24+18=\" },\n    @{ Old = \"67+4=\"; New = \"43+28=\" },\n    @{ Old = \"68+13=\"; New = \"40-11=\" },\n    @{ Old = \"59+19=\"; New = \"66-59=\" },\n    @{ Old = \"17+45=\"; New = \"71-6=\" },\n    @{ Old = \"74-38=\"; New = \"48+45=\" },\n    @{ Old = \"38+7=\"; New = \"2+59=\" },\n    @{ Old = \"83-44=\"; New = \"86+7=\" },\n    @{ Old = \"66-18=\"; New = \"58+29=\" },\n    @{ Old = \"87+5=\"; New = \"89+4=\" },\n    @{ Old = \"14+38=\"; New = \"40-37=\" },\n    @{ Old = \"12+49=\"; New = \"76+8=\" },\n    @{ Old = \"34-9=\"; New = \"56-29=\" },\n    @{ Old = \"8+84=\"; New = \"15+17=\" },\n    @{ Old = \"92-83=\"; New = \"31-6=\" },\n    @{ Old = \"39+49=\"; New = \"36+17=\" },\n    @{ Old = \"9+23=\"; New = \"63+28=\" },\n    @{ Old = \"67+8=\"; New = \"78+6=\" },\n    @{ Old = \"51-27=\"; New = \"69+14=\" },\n    @{ Old = \"23+59=\"; New = \"33-8=\" },\n    @{ Old = \"29+23=\"; New = \"3+88=\" },\n    @{ Old = \"19+68=\"; New = \"88+8=\" },\n    @{ Old = \"95-39=\"; New = \"19+56=\" },\n    @{ Old = \"18+79=\"; New = \"33+58=\" },\n    @{ Old = \"2+9=\"; New = \"37+54=\" },\n    @{ Old = \"49+6=\"; New = \"22-14=\" },\n    @{ Old = \"96-78=\"; New = \"81-49=\" },\n    @{ Old = \"28+8=\"; New = \"73-44=\" },\n    @{ Old = \"75-9=\"; New = \"82-3=\" },\n    @{ Old = \"30-16=\"; New = \"48+47=\" },\n    @{ Old = \"82-77=\"; New = \"90-16=\" },\n    @{ Old = \"33-14=\"; New = \"65+9=\" },\n    @{ Old = \"46-17=\"; New = \"6+19=\" },\n    @{ Old = \"91-86=\"; New = \"3+49=\" },\n    @{ Old = \"57+5=\"; New = \"88-69=\" },\n    @{ Old = \"15+68=\"; New = \"76-39=\" },\n    @{ Old = \"75-39=\"; New = \"45+8=\" },\n    @{ Old = \"90-12=\"; New = \"16+65=\" },\n    @{ Old = \"71-62=\"; New = \"3+79=\" },\n    @{ Old = \"61-58=\"; New = \"70-29=\" },\n    @{ Old = \"25-18=\"; New = \"36+5=\" }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$i = 0\nforeach ($cell in $t.Range.Cells) {\n    $pair = $pairs[$i]\n    $i++\n\n    $r = $cell.Range\n    # Drop the trailing paragraph mark + cell mark so only the visible\n    # text is compared/replaced; this keeps the run's rPr (font/size)\n    # intact instead of clearing formatting for the whole cell.\n    $r.MoveEnd(1, -1) | Out-Null\n\n    if ($r.Text -eq $pair.Old) {\n        $r.Text = $pair.New\n    }\n}\n"}
